# Automatische test-sync: 2025-08-19 19:17:50
# Append a new log row to "Logs" and bump the matching Dashboard count.

$wb = $excel.ActiveWorkbook

$wsLogs = $wb.Worksheets.Item("Logs")

# New row 3 of the Logs sheet (C3 and E3 stay empty, same as row 2).
$wsLogs.Range("A3").Value = "Interne taak"
$wsLogs.Range("B3").Value = "kwaliteit@testbedrijf123.nl"
$wsLogs.Range("D3").Value = "Intern verzoek / Actie voor medewerker"
$wsLogs.Range("F3").Value = "2025-08-19 19:17:12"
$wsLogs.Range("G3").Value = "Nee"
$wsLogs.Range("H3").Value = "Ja"
$wsLogs.Range("I3").Value = "Nee"
$wsLogs.Range("J3").Value = "Nee"

# Extend the existing conditional-formatting rules (previously scoped to row 2
# only) so they also cover the newly added row 3.
foreach ($col in "D","G","H","I","J") {
    $headCell = $wsLogs.Range("$($col)2")
    $newRange = $wsLogs.Range("$($col)2:$($col)3")
    $fc = $headCell.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Dashboard: bump the "Intern verzoek / Actie voor medewerker" count 1 -> 2.
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Range("B2").Value = 2
